# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they match the rest of the header row (bold font,
# centered/top aligned, thin box border), then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 72
    $ws.Cells.Item($row, 31).Value = 90
    $ws.Cells.Item($row, 32).Value = 0
}
